$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 28) mirroring the existing table rows.
$row = 28

$ws.Range("A$row").Value = "19/20"
$ws.Range("B$row").Value = "zat"

$ws.Range("C$row").Value = 43974.875
$ws.Range("C$row").NumberFormat = "DD/MM/YY\ HH:MM"

$ws.Range("D$row").Value = "Paul de Munnik"
$ws.Range("E$row").Value = "1.1 Poppodium Bolwerk"
$ws.Range("F$row").Value = "1 x voorst."
$ws.Range("G$row").Value = "BOL"
$ws.Range("H$row").Value = "Bolwerk Concert"
$ws.Range("I$row").Value = "Rock"
$ws.Range("J$row").Value = "BOL15+D18"
$ws.Range("O$row").Value = "Garderobetoeslag Bolwerk"
$ws.Range("P$row").Value = "Gepubliceerd"

$ws.Rows.Item($row).RowHeight = 13.8
